$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 26359.334
$ws.Range("J3").Value = 26359.334
$ws.Range("L3").Value = 26359.334
$ws.Range("N3").Value = -26587.334

$ws.Range("H26").Value = 15333.333
$ws.Range("J26").Value = 22000
$ws.Range("L26").Value = 22000
$ws.Range("N26").Value = -22688

$ws.Range("H28").Value = 653
$ws.Range("J28").Value = 1137.2
$ws.Range("L28").Value = 1137.2
$ws.Range("N28").Value = -2107.2

$ws.Range("H62").Value = 25008592
$ws.Range("I62").Value = 1291.2
$ws.Range("J62").Value = 100030500
$ws.Range("K62").Value = 1291.2
$ws.Range("L62").Value = 100030500
$ws.Range("M62").Value = -667.2
$ws.Range("N62").Value = -100031748

$ws.Range("H64").Value = 3445.5881
$ws.Range("I64").Value = 3070.4546
$ws.Range("J64").Value = 4133.3335
$ws.Range("K64").Value = 3070.4546
$ws.Range("L64").Value = 4133.3335
$ws.Range("M64").Value = -2822.4546
$ws.Range("N64").Value = -4629.3335

$ws.Range("H65").Value = 25008592
$ws.Range("I65").Value = 1291.2
$ws.Range("J65").Value = 100030500
$ws.Range("K65").Value = 6456
$ws.Range("L65").Value = 500152500
$ws.Range("M65").Value = -3336
$ws.Range("N65").Value = -500158740

$ws.Range("H67").Value = 3445.5881
$ws.Range("I67").Value = 3070.4546
$ws.Range("J67").Value = 4133.3335
$ws.Range("K67").Value = 3070.4546
$ws.Range("L67").Value = 4133.3335
$ws.Range("M67").Value = -2212.4546
$ws.Range("N67").Value = -5849.3335

$ws.Range("H100").Value = 7081
$ws.Range("I100").Value = 2493.3333
$ws.Range("J100").Value = 11668.667
$ws.Range("K100").Value = 2493.3333
$ws.Range("L100").Value = 11668.667
$ws.Range("M100").Value = -1952.3333
$ws.Range("N100").Value = -12750.667

$ws.Range("H102").Value = 26359.334
$ws.Range("J102").Value = 26359.334
$ws.Range("L102").Value = 26359.334
$ws.Range("N102").Value = -32849.334

$ws.Range("H105").Value = 30705
$ws.Range("J105").Value = 30705
$ws.Range("L105").Value = 30705
$ws.Range("N105").Value = -37693

$ws.Range("H138").Value = 23393982
$ws.Range("I138").Value = 90911320
$ws.Range("J138").Value = 7248533.5
$ws.Range("K138").Value = 272733960
$ws.Range("L138").Value = 21745600.5
$ws.Range("M138").Value = -272728820
$ws.Range("N138").Value = -21755880.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 23668.285
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H61").Value = 2053.92
$ws.Range("I61").Value = 1635.619
$ws.Range("J61").Value = 4250
$ws.Range("K61").Value = 1635.619
$ws.Range("L61").Value = 4250
$ws.Range("M61").Value = -1423.619
$ws.Range("N61").Value = -4674

$ws.Range("H136").Value = 2053.92
$ws.Range("I136").Value = 1635.619
$ws.Range("J136").Value = 4250
$ws.Range("K136").Value = 4906.857
$ws.Range("L136").Value = 12750
$ws.Range("M136").Value = -2356.857
$ws.Range("N136").Value = -17850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2842.7273
$ws.Range("I99").Value = 1721.25
$ws.Range("K99").Value = 1721.25
$ws.Range("M99").Value = -223.25

$ws.Range("H105").Value = 3863.3684
$ws.Range("I105").Value = 2955.5557
$ws.Range("J105").Value = 4680.4
$ws.Range("K105").Value = 2955.5557
$ws.Range("L105").Value = 4680.4
$ws.Range("M105").Value = -1208.5557
$ws.Range("N105").Value = -8174.4

$ws.Range("H107").Value = 8340.909
$ws.Range("I107").Value = 8805.556
$ws.Range("J107").Value = 6250
$ws.Range("K107").Value = 8805.556
$ws.Range("L107").Value = 6250
$ws.Range("M107").Value = -6885.556
$ws.Range("N107").Value = -10090

$ws.Range("H132").Value = 39910
$ws.Range("J132").Value = 39910
$ws.Range("L132").Value = 39910
$ws.Range("N132").Value = -50030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("N29").Value = 0

$ws.Range("H50").Value = 8342.299999999999
$ws.Range("J50").Value = 8824.777
$ws.Range("L50").Value = 8824.777
$ws.Range("N50").Value = -10074.777

$ws.Range("H68").Value = 35221.25
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 35221.25
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H97").Value = 14294.444
$ws.Range("J97").Value = 14294.444
$ws.Range("L97").Value = 14294.444
$ws.Range("N97").Value = -16276.444

$ws.Range("H105").Value = 2537.3684
$ws.Range("I105").Value = 1106.6666
$ws.Range("J105").Value = 4990
$ws.Range("K105").Value = 1106.6666
$ws.Range("L105").Value = 4990
$ws.Range("M105").Value = 640.3334
$ws.Range("N105").Value = -8484

$ws.Range("H107").Value = 1520.5883
$ws.Range("I107").Value = 1536.3636
$ws.Range("J107").Value = 1491.6666
$ws.Range("K107").Value = 1536.3636
$ws.Range("L107").Value = 1491.6666
$ws.Range("M107").Value = 383.6364000000001
$ws.Range("N107").Value = -5331.6666

$ws.Range("H109").Value = 20614
$ws.Range("J109").Value = 20614
$ws.Range("L109").Value = 20614
$ws.Range("N109").Value = -22694

$ws.Range("H123").Value = 50076.668
$ws.Range("J123").Value = 50076.668
$ws.Range("L123").Value = 50076.668
$ws.Range("N123").Value = -59876.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2000.8572
$ws.Range("J34").Value = 2000.8572
$ws.Range("L34").Value = 6002.571599999999
$ws.Range("N34").Value = -6170.571599999999

$ws.Range("H113").Value = 596.6087
$ws.Range("I113").Value = 530.8182
$ws.Range("J113").Value = 656.9167
$ws.Range("K113").Value = 1592.4546
$ws.Range("L113").Value = 1970.7501
$ws.Range("M113").Value = 577.5454
$ws.Range("N113").Value = -6310.7501

$ws.Range("H131").Value = 849.4815
$ws.Range("I131").Value = 375
$ws.Range("J131").Value = 867.7308
$ws.Range("K131").Value = 1125
$ws.Range("L131").Value = 2603.1924
$ws.Range("M131").Value = 3915
$ws.Range("N131").Value = -12683.1924

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1000
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -645

$ws.Range("H57").Value = 18983.691
$ws.Range("J57").Value = 18983.691
$ws.Range("L57").Value = 18983.691
$ws.Range("N57").Value = -20623.691

$ws.Range("H107").Value = 532.3
$ws.Range("I107").Value = 235.35715
$ws.Range("J107").Value = 1225.1666
$ws.Range("K107").Value = 235.35715
$ws.Range("L107").Value = 1225.1666
$ws.Range("M107").Value = 1684.64285
$ws.Range("N107").Value = -5065.1666

$ws.Range("H122").Value = 2221.15
$ws.Range("I122").Value = 2176.5
$ws.Range("J122").Value = 2399.75
$ws.Range("K122").Value = 6529.5
$ws.Range("L122").Value = 7199.25
$ws.Range("M122").Value = -4079.5
$ws.Range("N122").Value = -12099.25

$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -24900

$ws.Range("H124").Value = 49092
$ws.Range("J124").Value = 49092
$ws.Range("L124").Value = 49092
$ws.Range("N124").Value = -58912

$ws.Range("H133").Value = 39563
$ws.Range("J133").Value = 39563
$ws.Range("L133").Value = 39563
$ws.Range("N133").Value = -49683

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("N5").Value = 0

$ws.Range("H46").Value = 1310.0454
$ws.Range("I46").Value = 1109
$ws.Range("J46").Value = 1511.091
$ws.Range("K46").Value = 1109
$ws.Range("L46").Value = 1511.091
$ws.Range("M46").Value = -921
$ws.Range("N46").Value = -1887.091

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0

$ws.Range("H132").Value = 2372.8928
$ws.Range("I132").Value = 2081.65
$ws.Range("J132").Value = 3101
$ws.Range("K132").Value = 6244.950000000001
$ws.Range("L132").Value = 9303
$ws.Range("M132").Value = -3714.950000000001
$ws.Range("N132").Value = -14363

$ws.Range("H133").Value = 59964.2
$ws.Range("J133").Value = 59964.2
$ws.Range("L133").Value = 59964.2
$ws.Range("N133").Value = -65024.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 17944.25
$ws.Range("J21").Value = 17944.25
$ws.Range("L21").Value = 17944.25
$ws.Range("N21").Value = -18414.25

$ws.Range("H35").Value = 17944.25
$ws.Range("J35").Value = 17944.25
$ws.Range("L35").Value = 17944.25
$ws.Range("N35").Value = -18524.25

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0

$ws.Range("H109").Value = 29638.5
$ws.Range("J109").Value = 29638.5
$ws.Range("L109").Value = 29638.5
$ws.Range("N109").Value = -32412.5

$ws.Range("H132").Value = 4782.1763
$ws.Range("I132").Value = 4739.9
$ws.Range("J132").Value = 4842.5713
$ws.Range("K132").Value = 14219.7
$ws.Range("L132").Value = 14527.7139
$ws.Range("M132").Value = -11689.7
$ws.Range("N132").Value = -19587.7139

$ws.Range("H136").Value = 2298
$ws.Range("I136").Value = 1816.1177
$ws.Range("J136").Value = 3663.3333
$ws.Range("K136").Value = 5448.3531
$ws.Range("L136").Value = 10989.9999
$ws.Range("M136").Value = -2898.3531
$ws.Range("N136").Value = -16089.9999
